$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 3, 4, 7, 8
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -6
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = 2
